$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -12.992
$ws.Range("C21").Value = -12.499
$ws.Range("C23").Value = -12.594
$ws.Range("C25").Value = -12.37
$ws.Range("D27").Value = -8.484
$ws.Range("D31").Value = -8.312000000000001
$ws.Range("D39").Value = -7.747
$ws.Range("D48").Value = -7.475
$ws.Range("D51").Value = -8.434000000000001
$ws.Range("D52").Value = -7.568
$ws.Range("C53").Value = -11.193
$ws.Range("D55").Value = -8.065000000000001
$ws.Range("D56").Value = -8.288
$ws.Range("C57").Value = -13.565
$ws.Range("D57").Value = -8.559000000000001
$ws.Range("C59").Value = -13.077
$ws.Range("C69").Value = -10.676
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("C79").Value = -12.013
$ws.Range("C83").Value = -13.169
$ws.Range("D89").Value = -6.702
$ws.Range("D90").Value = -7.602000000000001
$ws.Range("C93").Value = -11.511

$wb.Save()
